# Skad_HRMS TestData.xlsx - "changes made in user ,base ,departments"
#
# Adds a new "VendorAdministration" sheet (vendor / company master data)
# after the existing "Data" and "ProjectData" sheets, populates its header
# row + one data row, wires up a mailto hyperlink on the email cell, and
# makes the new sheet the active / selected tab - mirroring what happens
# when a user adds & fills in a brand-new worksheet in the Excel UI.

$wb = $excel.ActiveWorkbook

# --- add the new sheet at the end of the tab strip -------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "VendorAdministration"

# --- header row --------------------------------------------------------
$ws.Range("A1").Value = "CompanyName"
$ws.Range("B1").Value = "Mobile"
$ws.Range("C1").Value = "Email"
$ws.Range("D1").Value = "Country"
$ws.Range("E1").Value = "State"
$ws.Range("F1").Value = "City"
$ws.Range("G1").Value = "PinCode"
$ws.Range("H1").Value = "PermanentAddress"

# --- first data row ------------------------------------------------------
$ws.Range("A2").Value = "Dolphin "
$ws.Range("B2").Value = 9159268812
$ws.Range("C2").Value = "dolphin@gmail.com"
$ws.Range("D2").Value = "India"
$ws.Range("E2").Value = "Tamil Nadu"
$ws.Range("F2").Value = "Thanjavur"
$ws.Range("G2").Value = 614625
$ws.Range("H2").Value = "3 North street Orathanadu"

# --- hyperlink the email address, like the existing sheets do ----------
$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:dolphin@gmail.com", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "dolphin@gmail.com") | Out-Null

# --- column widths (approximate the source layout) ----------------------
$ws.Columns.Item(1).ColumnWidth = 19.27
$ws.Columns.Item(2).ColumnWidth = 21.45
$ws.Columns.Item(3).ColumnWidth = 16.91
$ws.Columns.Item(4).ColumnWidth = 13.45
$ws.Columns.Item(5).ColumnWidth = 18.55
$ws.Columns.Item(6).ColumnWidth = 14.82
$ws.Columns.Item(7).ColumnWidth = 17.64
$ws.Columns.Item(8).ColumnWidth = 25.36

# --- make the new sheet the active / selected one -----------------------
$ws.Activate()
$ws.Range("H2").Select()
